$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.28083642826960842
$ws.Range("C2").Value = 0.86631504108253032
$ws.Range("B3").Value = 0.20315860626223561
$ws.Range("C3").ClearContents()

$ws.Range("A3:B3").Copy()
$ws.Range("A4:B5").PasteSpecial(-4122)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 0.19725091881722531

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0.18506908773346101

$excel.CutCopyMode = $false
